$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Transactions")
$cell = $ws.Range("B21")

# The new value is a purely numeric-looking string ("3397079808"), and the
# existing column B entries are likewise numeric-looking strings stored as
# text (shared strings), not numbers. Force text entry by formatting the
# cell as Text before assigning the value (otherwise Excel auto-converts a
# numeric-looking value to a real number), then restore the cell's style
# to the workbook's default "Normal" style so no stray per-cell formatting
# is left behind - matching the original sheet's unstyled text cells.
$cell.NumberFormat = "@"
$cell.Value = "3397079808"
$cell.Style = "Normal"
